$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "PlayFixedCutscene" right after the existing
# NpcMenuFunctionType / PlayCutscene-PlaySlimeMinigame rows, pushing the
# rest of the table down by one.
$ws.Rows(8).Insert()
$ws.Rows(8).RowHeight = 15.75

# The blank row created above inherits formatting from row 7; fill in the
# row that used to be there (PlaySlimeMinigame, now Order = 3) ...
$ws.Range("A8").Value = "NpcMenuFunctionType"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "PlaySlimeMinigame"

# ... and put the brand-new enum value where the old row used to be.
$ws.Range("C7").Value = "PlayFixedCutscene"

# Add a new trailing (mostly empty) row 22, matching the formatting-only
# placeholder left at the bottom of the table, then drop the extra
# formatted-but-empty cells it creates in columns B/C.
$ws.Rows(22).Insert()
$ws.Rows(22).RowHeight = 15.75
$ws.Range("B22:C22").Clear()

# Restore the active selection to the cell the author ended up editing.
$ws.Range("C8").Select()
